$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7828204035758972
$ws.Range("B1").Value = 2.239025115966797
$ws.Range("D1").Value = 0.9308872818946838
$ws.Range("E1").Value = 0.636476993560791
